$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: introduce the 5 brand-new shared strings in the exact order
#     they were first typed by the user (this controls sharedStrings.xml order) ---
$ws.Range("D54").Value = ' '
$ws.Range("N11").Value = 'xrxx'
$ws.Range("P5").Value = 'xxx?'
$ws.Range("AA8").Value = 'yry'
$ws.Range("N14").Value = '53mm'

# --- Step 2: fill in the remaining new cells for rows 5-12 (extra condition columns) ---
# Row 5
$ws.Range("N5").Value = 'y'
$ws.Range("O5").Value = 'y'
$ws.Range("R5").Value = 'y'
$ws.Range("S5").Value = 'y'
$ws.Range("T5").Value = 'y'
$ws.Range("V5").Value = 'y'
$ws.Range("W5").Value = 'y'
$ws.Range("X5").Value = 'y'
$ws.Range("Z5").Value = 'y'
$ws.Range("AA5").Value = 'y'
$ws.Range("AB5").Value = 'y'
$ws.Range("AD5").Value = 'y'
$ws.Range("AE5").Value = 'y'
$ws.Range("AF5").Value = 'y'
$ws.Range("AH5").Value = 'y'
$ws.Range("AI5").Value = 'y'
$ws.Range("AJ5").Value = 'y'
$ws.Range("AL5").Value = 'y'
$ws.Range("AM5").Value = 'y'
$ws.Range("AN5").Value = 'y'
# Row 6
$ws.Range("N6").Value = 'x'
$ws.Range("O6").Value = 'y'
$ws.Range("P6").Value = 'y'
$ws.Range("R6").Value = 'y'
$ws.Range("S6").Value = 'x'
$ws.Range("T6").Value = 'y'
$ws.Range("V6").Value = 'x'
$ws.Range("W6").Value = 'y'
$ws.Range("X6").Value = 'y'
$ws.Range("Z6").Value = 'y'
$ws.Range("AA6").Value = 'y'
$ws.Range("AB6").Value = 'y'
$ws.Range("AD6").Value = 'y'
$ws.Range("AE6").Value = 'y'
$ws.Range("AF6").Value = 'y'
$ws.Range("AH6").Value = 'y'
$ws.Range("AI6").Value = 'y'
$ws.Range("AJ6").Value = 'y'
$ws.Range("AL6").Value = 'y'
$ws.Range("AM6").Value = 'y'
$ws.Range("AN6").Value = 'y'
# Row 7
$ws.Range("N7").Value = 'y'
$ws.Range("O7").Value = 'y'
$ws.Range("P7").Value = 'y'
$ws.Range("R7").Value = 'xry'
$ws.Range("S7").Value = 'y'
$ws.Range("T7").Value = 'y'
$ws.Range("V7").Value = 'y'
$ws.Range("W7").Value = 'y'
$ws.Range("X7").Value = 'xrx'
$ws.Range("Z7").Value = 'y'
$ws.Range("AA7").Value = 'y'
$ws.Range("AB7").Value = 'y'
$ws.Range("AD7").Value = 'y'
$ws.Range("AE7").Value = 'y'
$ws.Range("AF7").Value = 'y'
$ws.Range("AH7").Value = 'y'
$ws.Range("AI7").Value = 'y'
$ws.Range("AJ7").Value = 'y'
$ws.Range("AL7").Value = 'y'
$ws.Range("AM7").Value = 'y'
$ws.Range("AN7").Value = 'y'
# Row 8
$ws.Range("N8").Value = 'y'
$ws.Range("O8").Value = 'y'
$ws.Range("P8").Value = 'y'
$ws.Range("R8").Value = 'y'
$ws.Range("S8").Value = 'y'
$ws.Range("T8").Value = 'y'
$ws.Range("V8").Value = 'y'
$ws.Range("W8").Value = 'y'
$ws.Range("X8").Value = 'y'
$ws.Range("Z8").Value = 'o'
$ws.Range("AB8").Value = 'y'
$ws.Range("AD8").Value = 'y'
$ws.Range("AE8").Value = 'y'
$ws.Range("AF8").Value = 'y'
$ws.Range("AH8").Value = 'y'
$ws.Range("AI8").Value = 'y'
$ws.Range("AJ8").Value = 'y'
$ws.Range("AL8").Value = 'y'
$ws.Range("AM8").Value = 'x'
$ws.Range("AN8").Value = 'o'
# Row 9
$ws.Range("N9").Value = 'y'
$ws.Range("O9").Value = 'y'
$ws.Range("P9").Value = 'y'
$ws.Range("R9").Value = 'y'
$ws.Range("S9").Value = 'y'
$ws.Range("T9").Value = 'y'
$ws.Range("V9").Value = 'y'
$ws.Range("W9").Value = 'y'
$ws.Range("X9").Value = 'y'
$ws.Range("Z9").Value = 'y'
$ws.Range("AA9").Value = 'y'
$ws.Range("AB9").Value = 'y'
$ws.Range("AD9").Value = 'y'
$ws.Range("AE9").Value = 'y'
$ws.Range("AF9").Value = 'y'
$ws.Range("AH9").Value = 'xry'
$ws.Range("AI9").Value = 'y'
$ws.Range("AJ9").Value = 'y'
$ws.Range("AL9").Value = 'xry'
$ws.Range("AM9").Value = 'y'
$ws.Range("AN9").Value = 'y'
# Row 10
$ws.Range("N10").Value = 'y'
$ws.Range("O10").Value = 'y'
$ws.Range("P10").Value = 'x'
$ws.Range("R10").Value = 'y'
$ws.Range("S10").Value = 'x'
$ws.Range("T10").Value = 'y'
$ws.Range("V10").Value = 'x'
$ws.Range("W10").Value = 'y'
$ws.Range("X10").Value = 'o'
$ws.Range("Z10").Value = 'y'
$ws.Range("AA10").Value = 'y'
$ws.Range("AB10").Value = 'x'
$ws.Range("AD10").Value = 'y'
$ws.Range("AE10").Value = 'y'
$ws.Range("AF10").Value = 'y'
$ws.Range("AH10").Value = 'y'
$ws.Range("AI10").Value = 'y'
$ws.Range("AJ10").Value = 'y'
$ws.Range("AL10").Value = 'o'
$ws.Range("AM10").Value = 'y'
$ws.Range("AN10").Value = 'y'
# Row 11
$ws.Range("O11").Value = 'y'
$ws.Range("P11").Value = 'x'
$ws.Range("R11").Value = 'x'
$ws.Range("S11").Value = 'xxry'
$ws.Range("T11").Value = 'y'
$ws.Range("V11").Value = 'o'
$ws.Range("W11").Value = 'y'
$ws.Range("X11").Value = 'y'
$ws.Range("Z11").Value = 'y'
$ws.Range("AA11").Value = 'y'
$ws.Range("AB11").Value = 'y'
$ws.Range("AD11").Value = 'y'
$ws.Range("AE11").Value = 'y'
$ws.Range("AF11").Value = 'y'
$ws.Range("AH11").Value = 'y'
$ws.Range("AI11").Value = 'y'
$ws.Range("AJ11").Value = 'y'
$ws.Range("AL11").Value = 'y'
$ws.Range("AM11").Value = 'y'
$ws.Range("AN11").Value = 'y'
# Row 12
$ws.Range("N12").Value = 'o'
$ws.Range("O12").Value = 'xry'
$ws.Range("P12").Value = 'y'
$ws.Range("R12").Value = 'o'
$ws.Range("S12").Value = 'y'
$ws.Range("T12").Value = 'x'
$ws.Range("V12").Value = 'y'
$ws.Range("W12").Value = 'o'
$ws.Range("X12").Value = 'o'
$ws.Range("Z12").Value = 'o'
$ws.Range("AA12").Value = 'x'
$ws.Range("AB12").Value = 'y'
$ws.Range("AD12").Value = 'y'
$ws.Range("AE12").Value = 'o'
$ws.Range("AF12").Value = 'y'
$ws.Range("AH12").Value = 'xry'
$ws.Range("AI12").Value = 'y'
$ws.Range("AJ12").Value = 'y'
$ws.Range("AL12").Value = 'o'
$ws.Range("AM12").Value = 'y'
$ws.Range("AN12").Value = 'o'

# --- Step 3: fill in the remaining new cells for rows 16-23 (extra condition columns) ---
# Row 16
$ws.Range("Z16").Value = 'y'
$ws.Range("AA16").Value = 'y'
$ws.Range("AB16").Value = 'y'
$ws.Range("AD16").Value = 'y'
$ws.Range("AE16").Value = 'y'
$ws.Range("AF16").Value = 'y'
$ws.Range("AH16").Value = 'y'
$ws.Range("AI16").Value = 'y'
$ws.Range("AJ16").Value = 'y'
# Row 17
$ws.Range("Z17").Value = 'y'
$ws.Range("AA17").Value = 'y'
$ws.Range("AB17").Value = 'y'
$ws.Range("AD17").Value = 'y'
$ws.Range("AE17").Value = 'y'
$ws.Range("AF17").Value = 'y'
$ws.Range("AH17").Value = 'y'
$ws.Range("AI17").Value = 'y'
$ws.Range("AJ17").Value = 'y'
# Row 18
$ws.Range("Z18").Value = 'y'
$ws.Range("AA18").Value = 'y'
$ws.Range("AB18").Value = 'y'
$ws.Range("AD18").Value = 'y'
$ws.Range("AE18").Value = 'y'
$ws.Range("AF18").Value = 'o'
$ws.Range("AH18").Value = 'y'
$ws.Range("AI18").Value = 'y'
$ws.Range("AJ18").Value = 'y'
# Row 19
$ws.Range("Z19").Value = 'y'
$ws.Range("AA19").Value = 'y'
$ws.Range("AB19").Value = 'y'
$ws.Range("AD19").Value = 'y'
$ws.Range("AE19").Value = 'y'
$ws.Range("AF19").Value = 'y'
$ws.Range("AH19").Value = 'y'
$ws.Range("AI19").Value = 'y'
$ws.Range("AJ19").Value = 'y'
# Row 20
$ws.Range("Z20").Value = 'y'
$ws.Range("AA20").Value = 'y'
$ws.Range("AB20").Value = 'y'
$ws.Range("AD20").Value = 'y'
$ws.Range("AE20").Value = 'y'
$ws.Range("AF20").Value = 'y'
$ws.Range("AH20").Value = 'y'
$ws.Range("AI20").Value = 'y'
$ws.Range("AJ20").Value = 'y'
# Row 21
$ws.Range("Z21").Value = 'y'
$ws.Range("AA21").Value = 'y'
$ws.Range("AB21").Value = 'y'
$ws.Range("AD21").Value = 'y'
$ws.Range("AE21").Value = 'y'
$ws.Range("AF21").Value = 'y'
$ws.Range("AH21").Value = 'y'
$ws.Range("AI21").Value = 'y'
$ws.Range("AJ21").Value = 'o'
# Row 22
$ws.Range("Z22").Value = 'y'
$ws.Range("AA22").Value = 'y'
$ws.Range("AB22").Value = 'y'
$ws.Range("AD22").Value = 'x'
$ws.Range("AE22").Value = 'y'
$ws.Range("AF22").Value = 'y'
$ws.Range("AH22").Value = 'y'
$ws.Range("AI22").Value = 'o'
$ws.Range("AJ22").Value = 'x'
# Row 23
$ws.Range("Z23").Value = 'y'
$ws.Range("AA23").Value = 'y'
$ws.Range("AB23").Value = 'y'
$ws.Range("AD23").Value = 'o'
$ws.Range("AE23").Value = 'y'
$ws.Range("AF23").Value = 'x'
$ws.Range("AH23").Value = 'y'
$ws.Range("AI23").Value = 'y'
$ws.Range("AJ23").Value = 'y'

# --- Step 4: update selection to match the saved view ---
$null = $ws.Range("N21").Select()
